# Generate Report for Archive
# Update the localization status report: files that have moved from
# "Ready for handoff" into the translation pipeline are now shown as
# "In Translation" on the Overview sheet as well as on each per-locale
# status sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B3").Value = "In Translation"
$ws.Range("C3").Value = "In Translation"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"

# --- zh-cn sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"

# --- de-de sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"
